# Update column G ("K") values for rows 2-19 per regenerated save_data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 0
    3  = 1
    4  = 2
    5  = 1
    6  = 1
    7  = 1
    8  = 1
    9  = 1
    11 = 0
    12 = 2
    13 = 3
    14 = 0
    15 = 0
    16 = 1
    17 = 1
    18 = 0
    19 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
